# RTM.xlsx update
# 1) update CRS version (V1.2 line of requirements bumped to V1.1)
# 2) edit RTM (HSI requirement references renumbered)
# 3) bump CR_overall_04 / CRS_overall_004 to V1.1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: CRS requirement version bump
$ws.Range("B2").Value = "Req_1ST123_CRS_overall_002-V1.1"

# Row 3: HSI requirement id updated
$ws.Range("C3").Value = "Req_1ST123_HSI_overall_002-V1.0"

# Row 4: HSI requirement ids updated (two-line cell)
$ws.Range("C4").Value = "Req_1ST123_HSI_overall_003-V1.0" + [char]10 + "Req_1ST123_HSI_overall_004-V1.0"

# Row 5: CR / CRS requirement versions bumped, HSI requirement id updated
$ws.Range("A5").Value = "Req_1ST123_CR_overall_04-V1.1"
$ws.Range("B5").Value = "Req_1ST123_CRS_overall_004-V1.1"
$ws.Range("C5").Value = "Req_1ST123_HSI_overall_005-V1.0"

# Leave the active selection on C5, matching the last-edited cell
$ws.Range("C5").Select()
